$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-0.0405101384793125, -0.0643051395617496, -0.0167151373968755),
    @(-0.0395525406328933, -0.0736720078149598, -0.00543307345082682),
    @(0.0252544925943794, -0.0598679690449487, 0.110376954233708),
    @(-0.0765248032469395, -0.148482662821163, -0.004566943672716),
    @(-0.171200400434677, -0.252797961937622, -0.0896028389317323),
    @(-0.0937905063089196, -0.197428467420491, 0.00984745480265203),
    @(0.00503032067922172, -0.0896479045450339, 0.0997085459034773),
    @(0.0434352917667326, -0.0399376672538099, 0.126808250787275),
    @(0.0845158701934025, -0.030081917545833, 0.199113657932638),
    @(-0.00527322075229465, -0.0576177795886535, 0.0470713380840642),
    @(-0.0530237098185993, -0.0941783085364134, -0.0118691111007853)
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
    $ws.Cells.Item($row, 4).Value = $values[$i][2]
}
